$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 24.06000000000032
$ws.Range("H2").Value = 0.01154131105659129
$ws.Range("I2").Value = 0.01154131105659129
$ws.Range("L2").Value = 34.40258958900493
$ws.Range("M2").Value = "[7.045151402152655, 61.76002777585721]"
$ws.Range("N2").Value = 0.01487262645182708
$ws.Range("O2").Value = 0.01487262645182708
$ws.Range("P2").Value = 1.742184514603349
$ws.Range("Q2").Value = "[0.6603948521059628, 2.8239741771007356]"
$ws.Range("R2").Value = 0.002227418392017011
$ws.Range("S2").Value = 0.002227418392017011
$ws.Range("T2").Value = 63.74209389249454
$ws.Range("U2").Value = "[48.12771180716884, 79.35647597782022]"
$ws.Range("V2").Value = [double]"1.64118052481399e-10"
$ws.Range("W2").Value = [double]"1.64118052481399e-10"
$ws.Range("X2").Value = 17.38870870870894
$ws.Range("Y2").Value = 13.24624624624642
$ws.Range("Z2").Value = 21.53117117117146

# Row 3
$ws.Range("F3").Value = 24.06000000000032
$ws.Range("H3").Value = [double]"2.980690203158165e-05"
$ws.Range("I3").Value = [double]"2.980690203158165e-05"
$ws.Range("L3").Value = 50.9131575289626
$ws.Range("M3").Value = "[25.743428144360763, 76.08288691356444]"
$ws.Range("N3").Value = 0.000185079666479071
$ws.Range("O3").Value = 0.000185079666479071
$ws.Range("P3").Value = 1.880552959806503
$ws.Range("Q3").Value = "[1.3145002294299646, 2.446605690183042]"
$ws.Range("R3").Value = [double]"2.928057796225403e-08"
$ws.Range("S3").Value = [double]"2.928057796225403e-08"
$ws.Range("T3").Value = 55.5529124472872
$ws.Range("U3").Value = "[41.47637531924923, 69.62944957532517]"
$ws.Range("V3").Value = [double]"4.094884431538048e-10"
$ws.Range("W3").Value = [double]"4.094884431538048e-10"
$ws.Range("X3").Value = 16.85885885885908
$ws.Range("Y3").Value = 14.69129129129149
$ws.Range("Z3").Value = 19.02642642642668

# Row 4
$ws.Range("F4").Value = 24.06000000000032
$ws.Range("H4").Value = [double]"1.443281539614816e-06"
$ws.Range("I4").Value = [double]"1.443281539614816e-06"
$ws.Range("L4").Value = 72.14891112961305
$ws.Range("M4").Value = "[44.12964868326438, 100.16817357596173]"
$ws.Range("N4").Value = [double]"4.935592084365226e-06"
$ws.Range("O4").Value = [double]"4.935592084365226e-06"
$ws.Range("P4").Value = 1.465447624197041
$ws.Range("Q4").Value = "[1.0126054398958102, 1.9182898084982716]"
$ws.Range("R4").Value = [double]"5.302956984643004e-08"
$ws.Range("S4").Value = [double]"5.302956984643004e-08"
$ws.Range("T4").Value = 66.88581737867972
$ws.Range("U4").Value = "[49.9425769093872, 83.82905784797224]"
$ws.Range("V4").Value = [double]"4.063178682400803e-10"
$ws.Range("W4").Value = [double]"4.063178682400803e-10"
$ws.Range("X4").Value = 18.44840840840865
$ws.Range("Y4").Value = 16.71435435435458
$ws.Range("Z4").Value = 20.18246246246273

# Row 5
$ws.Range("F5").Value = 24.06000000000032
$ws.Range("H5").Value = 0.0001481038694504111
$ws.Range("I5").Value = 0.0001481038694504111
$ws.Range("L5").Value = 45.7419787446825
$ws.Range("M5").Value = "[18.671385928420918, 72.81257156094408]"
$ws.Range("N5").Value = 0.001408302070097456
$ws.Range("O5").Value = 0.001408302070097456
$ws.Range("P5").Value = 2.081816152829272
$ws.Range("Q5").Value = "[1.4906055233248852, 2.6730267823336584]"
$ws.Range("R5").Value = [double]"7.438825111449887e-09"
$ws.Range("S5").Value = [double]"7.438825111449887e-09"
$ws.Range("T5").Value = 70.51444883764761
$ws.Range("U5").Value = "[56.49348003965956, 84.53541763563565]"
$ws.Range("V5").Value = [double]"3.47277762102749e-13"
$ws.Range("W5").Value = [double]"3.47277762102749e-13"
$ws.Range("X5").Value = 16.08816816816839
$ws.Range("Y5").Value = 13.82426426426445
$ws.Range("Z5").Value = 18.35207207207232

# Row 6
$ws.Range("F6").Value = 24.06000000000032
$ws.Range("H6").Value = 0.002403776890409937
$ws.Range("I6").Value = 0.002403776890409937
$ws.Range("L6").Value = 40.07718086422584
$ws.Range("M6").Value = "[11.739130519945476, 68.4152312085062]"
$ws.Range("N6").Value = 0.006600945703075567
$ws.Range("O6").Value = 0.006600945703075567
$ws.Range("P6").Value = 1.956026657190042
$ws.Range("Q6").Value = "[1.1761317842268095, 2.7359215301532736]"
$ws.Range("R6").Value = [double]"7.747936831448143e-06"
$ws.Range("S6").Value = [double]"7.747936831448143e-06"
$ws.Range("T6").Value = 60.65681322393115
$ws.Range("U6").Value = "[45.296492791080425, 76.01713365678187]"
$ws.Range("V6").Value = [double]"4.027951305829447e-10"
$ws.Range("W6").Value = [double]"4.027951305829447e-10"
$ws.Range("X6").Value = 16.56984984985007
$ws.Range("Y6").Value = 13.58342342342361
$ws.Range("Z6").Value = 19.55627627627654

# Row 7
$ws.Range("F7").Value = 24.06000000000032
$ws.Range("H7").Value = 0.008876424348425949
$ws.Range("I7").Value = 0.008876424348425949
$ws.Range("L7").Value = 38.53085531658884
$ws.Range("M7").Value = "[9.436609599645351, 67.62510103353233]"
$ws.Range("N7").Value = 0.01058696106330403
$ws.Range("O7").Value = 0.01058696106330403
$ws.Range("P7").Value = 1.767342413731195
$ws.Range("Q7").Value = "[0.7484474990534249, 2.7862373284089657]"
$ws.Range("R7").Value = 0.00108164850391157
$ws.Range("S7").Value = 0.00108164850391157
$ws.Range("T7").Value = 63.19601008381183
$ws.Range("U7").Value = "[46.258352708982414, 80.13366745864124]"
$ws.Range("V7").Value = [double]"1.768253321543511e-09"
$ws.Range("W7").Value = [double]"1.768253321543511e-09"
$ws.Range("X7").Value = 17.2923723723726
$ws.Range("Y7").Value = 13.39075075075093
$ws.Range("Z7").Value = 21.19399399399428

# Row 8
$ws.Range("F8").Value = 24.06000000000032
$ws.Range("H8").Value = 0.0002555789258631425
$ws.Range("I8").Value = 0.0002555789258631425
$ws.Range("L8").Value = 52.45079795371622
$ws.Range("M8").Value = "[22.696670283644707, 82.20492562378773]"
$ws.Range("N8").Value = 0.0009145416384501726
$ws.Range("O8").Value = 0.0009145416384501726
$ws.Range("P8").Value = 1.86797401024258
$ws.Range("Q8").Value = "[1.2012896833546574, 2.5346583371305034]"
$ws.Range("R8").Value = [double]"1.053617329382917e-06"
$ws.Range("S8").Value = [double]"1.053617329382917e-06"
$ws.Range("T8").Value = 72.90769150794651
$ws.Range("U8").Value = "[56.19566543995894, 89.61971757593409]"
$ws.Range("V8").Value = [double]"2.541522547971908e-11"
$ws.Range("W8").Value = [double]"2.541522547971908e-11"
$ws.Range("X8").Value = 16.90702702702725
$ws.Range("Y8").Value = 14.35411411411431
$ws.Range("Z8").Value = 19.4599399399402

# Row 9
$ws.Range("F9").Value = 24.06000000000032
$ws.Range("H9").Value = 0.008478360207880509
$ws.Range("I9").Value = 0.008478360207880509
$ws.Range("L9").Value = 36.86806211040219
$ws.Range("M9").Value = "[10.343068680225777, 63.39305554057861]"
$ws.Range("N9").Value = 0.007513049999859911
$ws.Range("O9").Value = 0.007513049999859911
$ws.Range("P9").Value = 1.490605523324887
$ws.Range("Q9").Value = "[0.4968685077749626, 2.4843425388748113]"
$ws.Range("R9").Value = 0.004142907235937354
$ws.Range("S9").Value = 0.004142907235937354
$ws.Range("T9").Value = 58.5557642117319
$ws.Range("U9").Value = "[42.4446932238383, 74.6668351996255]"
$ws.Range("V9").Value = [double]"3.421799510405776e-09"
$ws.Range("W9").Value = [double]"3.421799510405776e-09"
$ws.Range("X9").Value = 18.35207207207232
$ws.Range("Y9").Value = 14.54678678678698
$ws.Range("Z9").Value = 22.15735735735765

# Row 10
$ws.Range("F10").Value = 23.88000000000029
$ws.Range("H10").Value = [double]"3.026979884523939e-07"
$ws.Range("I10").Value = [double]"3.026979884523939e-07"
$ws.Range("L10").Value = 62.12718425500763
$ws.Range("M10").Value = "[40.4108260153195, 83.84354249469575]"
$ws.Range("N10").Value = [double]"7.033627806851683e-07"
$ws.Range("O10").Value = [double]"7.033627806851683e-07"
$ws.Range("P10").Value = 1.66671081721981
$ws.Range("Q10").Value = "[1.2641844311742716, 2.0692372032653488]"
$ws.Range("R10").Value = [double]"1.110180836150221e-10"
$ws.Range("S10").Value = [double]"1.110180836150221e-10"
$ws.Range("T10").Value = 64.84678819150049
$ws.Range("U10").Value = "[51.33799010754187, 78.35558627545912]"
$ws.Range("V10").Value = [double]"1.479705247220409e-12"
$ws.Range("W10").Value = [double]"1.479705247220409e-12"
$ws.Range("X10").Value = 17.54546546546568
$ws.Range("Y10").Value = 16.01561561561581
$ws.Range("Z10").Value = 19.07531531531555

# Row 11
$ws.Range("F11").Value = 23.88000000000029
$ws.Range("H11").Value = [double]"8.516838301619867e-05"
$ws.Range("I11").Value = [double]"8.516838301619867e-05"
$ws.Range("L11").Value = 56.67997414294188
$ws.Range("M11").Value = "[26.589172788786186, 86.77077549709757]"
$ws.Range("N11").Value = 0.0004401720840117029
$ws.Range("O11").Value = 0.0004401720840117029
$ws.Range("P11").Value = 2.03150035457358
$ws.Range("Q11").Value = "[1.440289725069194, 2.6227109840779654]"
$ws.Range("R11").Value = [double]"1.335536725655118e-08"
$ws.Range("S11").Value = [double]"1.335536725655118e-08"
$ws.Range("T11").Value = 69.34098858077625
$ws.Range("U11").Value = "[52.627932382301054, 86.05404477925146]"
$ws.Range("V11").Value = [double]"1.050348696907122e-10"
$ws.Range("W11").Value = [double]"1.050348696907122e-10"
$ws.Range("X11").Value = 16.15903903903924
$ws.Range("Y11").Value = 13.91207207207224
$ws.Range("Z11").Value = 18.40600600600623

# Row 12
$ws.Range("F12").Value = 23.88000000000029
$ws.Range("H12").Value = [double]"7.735110634787823e-06"
$ws.Range("I12").Value = [double]"7.735110634787823e-06"
$ws.Range("L12").Value = 60.17317521001116
$ws.Range("M12").Value = "[34.13379308744079, 86.21255733258153]"
$ws.Range("N12").Value = [double]"2.88045958924954e-05"
$ws.Range("O12").Value = [double]"2.88045958924954e-05"
$ws.Range("P12").Value = 1.742184514603348
$ws.Range("Q12").Value = "[1.2264475824825016, 2.257921446724195]"
$ws.Range("R12").Value = [double]"1.992856279997568e-08"
$ws.Range("S12").Value = [double]"1.992856279997568e-08"
$ws.Range("T12").Value = 66.07866741677475
$ws.Range("U12").Value = "[50.7116875583847, 81.44564727516479]"
$ws.Range("V12").Value = [double]"3.842548501609144e-11"
$ws.Range("W12").Value = [double]"3.842548501609144e-11"
$ws.Range("X12").Value = 17.25861861861883
$ws.Range("Y12").Value = 15.29849849849869
$ws.Range("Z12").Value = 19.21873873873898

# Row 13
$ws.Range("F13").Value = 23.88000000000029
$ws.Range("H13").Value = [double]"8.75058491495162e-05"
$ws.Range("I13").Value = [double]"8.75058491495162e-05"
$ws.Range("L13").Value = 44.61206575577513
$ws.Range("M13").Value = "[20.6573946650376, 68.56673684651265]"
$ws.Range("N13").Value = 0.0005014163589183163
$ws.Range("O13").Value = 0.0005014163589183163
$ws.Range("P13").Value = 1.905710858934349
$ws.Range("Q13").Value = "[1.3019212798660398, 2.5095004380026573]"
$ws.Range("R13").Value = [double]"9.200514750062894e-08"
$ws.Range("S13").Value = [double]"9.200514750062894e-08"
$ws.Range("T13").Value = 43.83313066997
$ws.Range("U13").Value = "[30.660360138930983, 57.00590120100901]"
$ws.Range("V13").Value = [double]"2.82234029480577e-08"
$ws.Range("W13").Value = [double]"2.82234029480577e-08"
$ws.Range("X13").Value = 16.63711711711732
$ws.Range("Y13").Value = 14.34234234234252
$ws.Range("Z13").Value = 18.93189189189213

# Row 14
$ws.Range("F14").Value = 23.88000000000029
$ws.Range("H14").Value = 0.001989435476664747
$ws.Range("I14").Value = 0.001989435476664747
$ws.Range("L14").Value = 43.63788963053233
$ws.Range("M14").Value = "[13.686010315436178, 73.58976894562848]"
$ws.Range("N14").Value = 0.005244394079811476
$ws.Range("O14").Value = 0.005244394079811476
$ws.Range("P14").Value = 1.918289808498272
$ws.Range("Q14").Value = "[1.1635528346628874, 2.6730267823336575]"
$ws.Range("R14").Value = [double]"6.179636126368848e-06"
$ws.Range("S14").Value = [double]"6.179636126368848e-06"
$ws.Range("T14").Value = 66.47958242216373
$ws.Range("U14").Value = "[50.101926769583926, 82.85723807474353]"
$ws.Range("V14").Value = [double]"1.916562464288063e-10"
$ws.Range("W14").Value = [double]"1.916562464288063e-10"
$ws.Range("X14").Value = 16.58930930930951
$ws.Range("Y14").Value = 13.72084084084101
$ws.Range("Z14").Value = 19.45777777777801
